$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Description" text for the opening remarks row: drop the
# "entry survey" mention.
$ws.Range("E2").Value = "Greetings from the FAMPS and FSN Chairs"

# Update the "Description" text for the Day 1 Wrap-up row: replace the
# exit-survey/networking blurb with the closing remarks/Day 2 preview text.
$ws.Range("E10").Value = "Closing from the FAMPS and FSN Chairs; Preview of Day 2"

# Move the active selection to E15, matching the saved view state.
$ws.Range("E15").Select()
